$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1744.4445
$ws.Range("I100").Value = 425
$ws.Range("J100").Value = 2800
$ws.Range("K100").Value = 425
$ws.Range("L100").Value = 2800
$ws.Range("M100").Value = 116
$ws.Range("N100").Value = -3882

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2238.09
$ws.Range("I32").Value = 2238.09
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2238.09
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1951.09
$ws.Range("N32").ClearContents()

$ws.Range("H45").Value = 83334130
$ws.Range("I45").Value = 111111740
$ws.Range("J45").Value = 1295
$ws.Range("K45").Value = 111111740
$ws.Range("L45").Value = 1295
$ws.Range("M45").Value = -111111363
$ws.Range("N45").Value = -2049

$ws.Range("H132").Value = 5180.661
$ws.Range("I132").Value = 5684.896
$ws.Range("J132").Value = 3451.8572
$ws.Range("K132").Value = 17054.688
$ws.Range("L132").Value = 10355.5716
$ws.Range("M132").Value = -14524.688
$ws.Range("N132").Value = -15415.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4901.5835
$ws.Range("I20").Value = 5564.375
$ws.Range("K20").Value = 5564.375
$ws.Range("M20").Value = -5317.375

$ws.Range("H86").Value = 2100.4194
$ws.Range("I86").Value = 1890.0526
$ws.Range("K86").Value = 1890.0526
$ws.Range("M86").Value = -767.0526

$ws.Range("H89").Value = 2100.4194
$ws.Range("I89").Value = 1890.0526
$ws.Range("K89").Value = 9450.262999999999
$ws.Range("M89").Value = -3834.262999999999

$ws.Range("H105").Value = 3413.0588
$ws.Range("I105").Value = 2504.4
$ws.Range("J105").Value = 3791.6667
$ws.Range("K105").Value = 2504.4
$ws.Range("L105").Value = 3791.6667
$ws.Range("M105").Value = -757.4000000000001
$ws.Range("N105").Value = -7285.6667

$ws.Range("H107").Value = 1392.7222
$ws.Range("I107").Value = 1093.5172
$ws.Range("J107").Value = 2632.2856
$ws.Range("K107").Value = 1093.5172
$ws.Range("L107").Value = 2632.2856
$ws.Range("M107").Value = 826.4828
$ws.Range("N107").Value = -6472.2856

$ws.Range("H134").Value = 4031.7292
$ws.Range("I134").Value = 4333.324
$ws.Range("K134").Value = 12999.972
$ws.Range("M134").Value = -10464.972

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3549452
$ws.Range("I31").Value = 2579.9429
$ws.Range("J31").Value = 13894495
$ws.Range("K31").Value = 2579.9429
$ws.Range("L31").Value = 13894495
$ws.Range("M31").Value = -2284.9429
$ws.Range("N31").Value = -13895085

$ws.Range("H34").Value = 3549452
$ws.Range("I34").Value = 2579.9429
$ws.Range("J34").Value = 13894495
$ws.Range("K34").Value = 2579.9429
$ws.Range("L34").Value = 13894495
$ws.Range("M34").Value = -2377.9429
$ws.Range("N34").Value = -13894899

$ws.Range("H58").Value = 1016.1778
$ws.Range("I58").Value = 1059.5264
$ws.Range("J58").Value = 780.8570999999999
$ws.Range("K58").Value = 1059.5264
$ws.Range("L58").Value = 780.8570999999999
$ws.Range("M58").Value = -856.5264
$ws.Range("N58").Value = -1186.8571

$ws.Range("H86").Value = 90910810
$ws.Range("I86").Value = 100001440
$ws.Range("K86").Value = 100001440
$ws.Range("M86").Value = -100000317

$ws.Range("H89").Value = 90910810
$ws.Range("I89").Value = 100001440
$ws.Range("K89").Value = 500007200
$ws.Range("M89").Value = -500001584

$ws.Range("H136").Value = 1016.1778
$ws.Range("I136").Value = 1059.5264
$ws.Range("J136").Value = 780.8570999999999
$ws.Range("K136").Value = 3178.5792
$ws.Range("L136").Value = 2342.5713
$ws.Range("M136").Value = -628.5792000000001
$ws.Range("N136").Value = -7442.5713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1259.2059
$ws.Range("I5").Value = 252.82353
$ws.Range("K5").Value = 758.47059
$ws.Range("M5").Value = -646.47059

$ws.Range("H18").Value = 418.85715
$ws.Range("I18").Value = 390
$ws.Range("J18").Value = 491
$ws.Range("K18").Value = 1170
$ws.Range("L18").Value = 1473
$ws.Range("M18").Value = -1001
$ws.Range("N18").Value = -1811

$ws.Range("H62").Value = 2692.5715
$ws.Range("I62").Value = 2319.8
$ws.Range("K62").Value = 6959.400000000001
$ws.Range("M62").Value = -6273.400000000001

$ws.Range("H65").Value = 2692.5715
$ws.Range("I65").Value = 2319.8
$ws.Range("K65").Value = 20878.2
$ws.Range("M65").Value = -17446.2

$ws.Range("H92").Value = 640.5
$ws.Range("I92").Value = 301
$ws.Range("J92").Value = 980
$ws.Range("K92").Value = 903
$ws.Range("L92").Value = 2940
$ws.Range("M92").Value = 345
$ws.Range("N92").Value = -5436

$ws.Range("H121").Value = 7168.9062
$ws.Range("I121").Value = 10360.6
$ws.Range("J121").Value = 5718.136
$ws.Range("K121").Value = 31081.8
$ws.Range("L121").Value = 17154.408
$ws.Range("M121").Value = -29771.8
$ws.Range("N121").Value = -19774.408

$ws.Range("H131").Value = 722.2347
$ws.Range("I131").Value = 280.5
$ws.Range("J131").Value = 821.625
$ws.Range("K131").Value = 841.5
$ws.Range("L131").Value = 2464.875
$ws.Range("M131").Value = 4198.5
$ws.Range("N131").Value = -12544.875

$ws.Range("H135").Value = 1259.2059
$ws.Range("I135").Value = 252.82353
$ws.Range("K135").Value = 2275.41177
$ws.Range("M135").Value = 259.5882299999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2022.1364
$ws.Range("I97").Value = 1948.3572
$ws.Range("J97").Value = 2151.25
$ws.Range("K97").Value = 1948.3572
$ws.Range("L97").Value = 2151.25
$ws.Range("M97").Value = -1452.3572
$ws.Range("N97").Value = -3143.25

$ws.Range("H123").Value = 28188.666
$ws.Range("J123").Value = 28188.666
$ws.Range("L123").Value = 28188.666
$ws.Range("N123").Value = -33088.666

$ws.Range("H132").Value = 4599.222
$ws.Range("I132").Value = 4562.879
$ws.Range("K132").Value = 13688.637
$ws.Range("M132").Value = -11158.637

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 10418198
$ws.Range("I68").Value = 1274
$ws.Range("J68").Value = 17858858
$ws.Range("K68").Value = 1274
$ws.Range("L68").Value = 17858858
$ws.Range("M68").Value = -525
$ws.Range("N68").Value = -17860356

$ws.Range("H71").Value = 10418198
$ws.Range("I71").Value = 1274
$ws.Range("J71").Value = 17858858
$ws.Range("K71").Value = 6370
$ws.Range("L71").Value = 89294290
$ws.Range("M71").Value = -2626
$ws.Range("N71").Value = -89301778

$ws.Range("H82").Value = 1683.7142
$ws.Range("I82").Value = 1627.2
$ws.Range("J82").Value = 1825
$ws.Range("K82").Value = 1627.2
$ws.Range("L82").Value = 1825
$ws.Range("M82").Value = -1266.2
$ws.Range("N82").Value = -2547

$ws.Range("H85").Value = 1683.7142
$ws.Range("I85").Value = 1627.2
$ws.Range("J85").Value = 1825
$ws.Range("K85").Value = 1627.2
$ws.Range("L85").Value = 1825
$ws.Range("M85").Value = -379.2
$ws.Range("N85").Value = -4321

$ws.Range("H100").Value = 47622724
$ws.Range("I100").Value = 166667680
$ws.Range("J100").Value = 4744.8
$ws.Range("K100").Value = 166667680
$ws.Range("L100").Value = 4744.8
$ws.Range("M100").Value = -166667139
$ws.Range("N100").Value = -5826.8

$ws.Range("H132").Value = 6776.795
$ws.Range("I132").Value = 8734.615
$ws.Range("K132").Value = 26203.845
$ws.Range("M132").Value = -23673.845

$ws.Range("H136").Value = 8253.105
$ws.Range("I136").Value = 12731.3
$ws.Range("K136").Value = 38193.89999999999
$ws.Range("M136").Value = -35643.89999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 13889214
$ws.Range("I107").Value = 5952723.5
$ws.Range("J107").Value = 27778074
$ws.Range("K107").Value = 17858170.5
$ws.Range("L107").Value = 83334222
$ws.Range("M107").Value = -17856250.5
$ws.Range("N107").Value = -83338062

Write-Output "Applied Garuda_Profits updates"
